$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.245.35'
$ws.Range("E2").Value = '  +4.47%  '

$ws.Range("D3").Value = '2.590.25'
$ws.Range("E3").Value = '  +2.31%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '''521.76'
$ws.Range("E5").Value = '  +1.85%  '

$ws.Range("D6").Value = '''140.79'
$ws.Range("E6").Value = '  +0.85%  '

$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("E8").Value = '  +1.79%  '

$ws.Range("D9").Value = '2.612.60'
$ws.Range("E9").Value = '  +2.94%  '

$ws.Range("D10").Value = '''6.51'
$ws.Range("E10").Value = '  -0.34%  '

$ws.Range("E11").Value = '  +2.21%  '

$ws.Range("E12").Value = '  +3.26%  '

$ws.Range("E13").Value = '  +2.69%  '

$ws.Range("D14").Value = '3.051.78'
$ws.Range("E14").Value = '  +2.26%  '

$ws.Range("D15").Value = '59.157.67'
$ws.Range("E15").Value = '  +4.25%  '

$ws.Range("D16").Value = '''20.54'
$ws.Range("E16").Value = '  +2.91%  '

$ws.Range("D17").Value = '2.607.41'
$ws.Range("E17").Value = '  +2.08%  '

$ws.Range("E18").Value = '  +0.65%  '

$ws.Range("D19").Value = '''339.81'
$ws.Range("E19").Value = '  +2.91%  '

$ws.Range("E20").Value = '  +1.86%  '

$ws.Range("D21").Value = '''10.21'
$ws.Range("E21").Value = '  +1.74%  '

$ws.Range("D22").Value = '''6.51'
$ws.Range("E22").Value = '  +6.28%  '

$ws.Range("E23").Value = '  -0.32%  '

$ws.Range("D24").Value = '''66.11'
$ws.Range("E24").Value = '  +1.80%  '

$ws.Range("E25").Value = '  +1.56%  '

$ws.Range("E26").Value = '  +1.96%  '

$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("D28").Value = '''7.15'
$ws.Range("E28").Value = '  +4.34%  '

$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("E30").Value = '  -2.41%  '

$ws.Range("E31").Value = '  -4.55%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '''18.83'
$ws.Range("E32").Value = '  +2.07%  '

$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '''1.57'
$ws.Range("E33").Value = '  +1.63%  '

$ws.Range("D34").Value = '''148.90'
$ws.Range("E34").Value = '  +0.56%  '

$ws.Range("E35").Value = '  +1.10%  '

$ws.Range("E36").Value = '  +0.80%  '

$ws.Range("D37").Value = '''36.35'
$ws.Range("E37").Value = '  +2.22%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '''1.46'
$ws.Range("E38").Value = '  +3.36%  '

$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").Value = '''0.836'
$ws.Range("E39").Value = '  +2.32%  '

$ws.Range("D40").Value = '''0.827'
$ws.Range("E40").Value = '  -1.82%  '

$ws.Range("E41").Value = '  +2.70%  '

$ws.Range("D42").Value = '''277.30'
$ws.Range("E42").Value = '  +6.51%  '

$ws.Range("D43").Value = '''0.998'
$ws.Range("E43").Value = '  -0.21%  '

$ws.Range("D44").Value = '''10.75'
$ws.Range("E44").Value = '  +1.30%  '

$ws.Range("E45").Value = '  +2.67%  '

$ws.Range("E46").Value = '  +0.50%  '

$ws.Range("E47").Value = '  +1.19%  '

$ws.Range("E48").Value = '  +0.53%  '

$ws.Range("D49").Value = '1.986.93'
$ws.Range("E49").Value = '  +1.26%  '

$ws.Range("E50").Value = '  +0.77%  '

$ws.Range("E51").Value = '  +0.38%  '
